# Refresh cryptos list: prices and 1h volume % changes updated,
# plus a couple of rank-table rows whose coin order was swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Force the literal text into the cell (no numeric/date auto-coercion),
    # then drop the temporary text format so the cell keeps its original style.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextCell 'D2' '25.824.25'
Set-TextCell 'E2' '  -0.07%  '
Set-TextCell 'D3' '1.635.71'
Set-TextCell 'E3' '  -0.07%  '
Set-TextCell 'E4' '  -0.03%  '
Set-TextCell 'D5' '215.40'
Set-TextCell 'E5' '  -0.29%  '
Set-TextCell 'D6' '0.505'
Set-TextCell 'E6' '  -0.36%  '
Set-TextCell 'E7' '  -0.06%  '
Set-TextCell 'D8' '0.258'
Set-TextCell 'E8' '  -0.05%  '
Set-TextCell 'E9' '  -0.03%  '
Set-TextCell 'D10' '19.90'
Set-TextCell 'E10' '  +1.72%  '
Set-TextCell 'D11' '0.0784'
Set-TextCell 'E11' '  +0.72%  '
Set-TextCell 'E12' '  -0.66%  '
Set-TextCell 'D13' '1.641.13'
Set-TextCell 'E13' '  +0.39%  '
Set-TextCell 'D14' '1.860.82'
Set-TextCell 'E14' '  -0.09%  '
Set-TextCell 'D15' '0.558'
Set-TextCell 'E15' '  -0.79%  '
Set-TextCell 'D16' '0.0₃0771'
Set-TextCell 'E16' '  +1.56%  '
Set-TextCell 'D17' '63.13'
Set-TextCell 'E17' '  -0.09%  '
Set-TextCell 'D18' '25.831.95'
Set-TextCell 'E18' '  -0.11%  '
Set-TextCell 'E19' '  -0.11%  '
Set-TextCell 'D20' '4.43'
Set-TextCell 'E20' '  +2.30%  '
Set-TextCell 'D21' '194.37'
Set-TextCell 'E21' '  -0.04%  '
Set-TextCell 'E22' '  +0.79%  '
Set-TextCell 'D23' '6.16'
Set-TextCell 'E23' '  +1.59%  '
Set-TextCell 'E24' '  -0.01%  '
Set-TextCell 'D25' '1.77'
Set-TextCell 'E25' '  -0.95%  '
Set-TextCell 'D26' '139.14'
Set-TextCell 'E26' '  -0.78%  '
Set-TextCell 'E27' '  -5.30%  '
Set-TextCell 'E28' '  +1.19%  '
Set-TextCell 'E29' '  +0.83%  '
Set-TextCell 'E30' '  +0.03%  '
Set-TextCell 'D31' '0.0497'
Set-TextCell 'E31' '  +1.72%  '
Set-TextCell 'E32' '  +1.05%  '
Set-TextCell 'E33' '  +1.13%  '
Set-TextCell 'E34' '  +2.38%  '
Set-TextCell 'D35' '2.38'
Set-TextCell 'E35' '  +0.24%  '
Set-TextCell 'D36' '0.902'
Set-TextCell 'E36' '  +0.28%  '
Set-TextCell 'E37' '  +0.16%  '
Set-TextCell 'E38' '  +0.67%  '
Set-TextCell 'D39' '1.111.31'
Set-TextCell 'E39' '  -1.65%  '
Set-TextCell 'D40' '0.0157'
Set-TextCell 'E40' '  +0.37%  '
Set-TextCell 'E41' '  +0.55%  '
Set-TextCell 'E42' '  +0.75%  '
Set-TextCell 'B43' 'Quant'
Set-TextCell 'C43' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D43' '99.34'
Set-TextCell 'E43' '  +1.95%  '
Set-TextCell 'B44' 'TrustWalletToken'
Set-TextCell 'C44' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D44' '0.803'
Set-TextCell 'E44' '  +0.32%  '
Set-TextCell 'E45' '  +0.23%  '
Set-TextCell 'B46' 'Aave'
Set-TextCell 'C46' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D46' '55.57'
Set-TextCell 'E46' '  +0.10%  '
Set-TextCell 'B47' 'SynthetixNetwork'
Set-TextCell 'C47' 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextCell 'D47' '2.50'
Set-TextCell 'E47' '  +11.57%  '
Set-TextCell 'E48' '  -5.76%  '
Set-TextCell 'D49' '7.70'
Set-TextCell 'E49' '  +0.16%  '
Set-TextCell 'E50' '  -0.43%  '
Set-TextCell 'D51' '1.01'
Set-TextCell 'E51' '  +0.20%  '
